$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-B64Value {
    param($cellRef, $b64)
    $bytes = [System.Convert]::FromBase64String($b64)
    $str = [System.Text.Encoding]::UTF8.GetString($bytes)
    $ws.Range($cellRef).Value = $str
}

# Header row: add entry_point/solutions/result/pass@1/pass@5/pass@10 columns
$ws.Range("E1").Value = "entry_point"
$ws.Range("F1").Value = "solutions"
$ws.Range("G1").Value = "result"
$ws.Range("H1").Value = "pass@1"
$ws.Range("I1").Value = "pass@5"
$ws.Range("J1").Value = "pass@10"

# Row 2
Set-B64Value "B2" "RGVmaW5lIGEgY2xhc3MgbmFtZWQgUm91dGVEb21haW4gdGhhdCBoYXMgYW4gYXRyaWJ1dHRlIGNhbGxlZCBmaWxlLiBUaGUgYXR0cmlidXRlIGNhbGxlZCBmaWxlIGlzIGVxdWFsIHRvIGFuIGlucHV0IGdpdmVuIHRvIHRoZSBjbGFzcy4gVGhlIGlucHV0IGdpdmVuIHRvIHRoZSBjbGFzcyBpcyBjYWxsZWQgZmlsZXBhdGgu"
Set-B64Value "D2" "Y2xhc3MgUm91dGVEb21haW4oKQogICAgICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6CiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGg="
$ws.Range("E2").Value = "test0"
Set-B64Value "F2" "WyIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcblxuICAgIGRlZiByZWFkX2ZpbGUoc2VsZik6XG4gICAgICAgIGZpbGUgPSBvcGVuKHNlbGYuZmlsZSlcbiAgICAgICAgZmlsZV9saXN0ID0gZmlsZS5yZWFkKCkuc3BsaXRsaW5lcygpXG4gICAgICAgIHJldHVybiBmaWxlX2xpc3RcblxuICAgIGRlZiByb3V0ZV9kb21haW4oc2VsZik6XG4gICAgICAgIGZpbGVfbGlzdCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgcm91dGVfZG9tYWluX2xpc3QgPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiBmaWxlX2xpc3Q6XG4gICAgICAgICAgICBsaW5lID0gbGluZS5zcGxpdCgnXHQnKVxuICAgICAgICAgICAgcm91dGVfZG9tYWluX2xpc3QuYXBwZW5kKGxpbmUpXG4gICAgICAgIHJldHVybiByb3V0ZV9kb21haW5fbGlzdFxuXG4gICAgZGVmIGRvbWFpbl9saXN0KHNlbGYpOlxuICAgICAgICByb3V0ZV9kb21haW5fbGlzdCA9IHNlbGYucm91dGVfZG9tYWluKClcbiAgICAgICAgZG9tYWluX2xpc3QgPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiByb3V0ZV9kb21haW5fbGlzdDpcbiAgICAgICAgICAgIGRvbWFpbl9saXN0LmFwcGVuZChsaW5lWzFdKVxuICAgICAgICByZXR1cm4gZG9tYWluX2xpc3RcblxuICAgIGRlZiByb3V0ZV9saXN0KHNlbGYpOlxuICAgICAgICByb3V0ZV9kb21haW5fbGlzdCA9IHNlbGYucm91dGVfZG9tYWluKClcbiAgICAgICAgcm91dGVfbGlzdCA9IFtdXG4gICAgICAgIGZvciBsaW5lIGluIHJvdXRlX2RvbWFpbl9saXN0OlxuICAgICAgICAgICAgcm91dGVfbGlzdC5hcHBlbmQobGluZVswXSlcbiAgICAgICAgcmV0dXJuIHJvdXRlX2xpc3RcblxuY2xhc3MgUm91dGUoKTpcbiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuXG4gICAgZGVmIHJlYWRfZmlsZShzZWxmKTpcbiAgICAgICAgZmlsZSA9IG9wZW4oc2VsZi5maWxlKVxuICAgICAgICBmaWxlX2xpc3QgPSBmaWxlLnJlYWQoKS5zcGxpdGxpbmVzKClcbiAgICAgICAgcmV0dXJuIGZpbGVfbGlzdFxuXG4gICAgZGVmIHJvdXRlX2xpc3Qoc2VsZik6XG4gICAgICAgIGZpbGVfbGlzdCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgcm91dGVfbGlzdCA9IFtdXG4gICAgICAgIGZvciBsaW5lIGluIGZpbGVfbGlzdDpcbiAgICAgICAgICAgIGxpbmUgPSBsaW5lLnNwbGl0KCdcdCcpXG4gICAgICAgICAgICByb3V0ZV9saXN0LmFwcGVuZChsaW5lKVxuICAgICAgICByZXR1cm4gcm91dGVfbGlzdFxuXG4gICAgZGVmIGRvbWFpbl9saXN0KHNlbGYpOlxuICAgICAgICByb3V0ZV9saXN0ID0gc2VsZi5yb3V0ZV9saXN0KClcbiAgICAgICAgZG9tYWluX2xpc3QgPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiByb3V0ZV9saXN0OlxuICAgICAgICAgICAgZG9tYWluX2xpc3QuYXBwZW5kKGxpbmVbMV0pXG4gICAgICAgIHJldHVybiBkb21haW5fbGlzdFxuXG4gICAgZGVmIHJvdXRlX2xpc3Qoc2VsZik6XG4gICAgICAgIHJvdXRlX2xpc3QgPSBzZWxmLnJvdXRlX2xpc3QoKVxuICAgICAgICByb3V0ZV9saXN0ID0gW11cbiAgICAgICAgZm9yIGxpbmUgaW4gcm91dGVfbGlzdDpcbiAgICAgICAgICAgIHJvdXRlX2xpc3QuYXBwZW5kKGxpbmVbMF0pXG4gICAgICAgIHJldHVybiByb3V0ZV9saXN0XG5cbiIsICcgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcblxuICAgIGRlZiBfX3N0cl9fKHNlbGYpOlxuICAgICAgICByZXR1cm4gIlJvdXRlRG9tYWluKHt9KSIuZm9ybWF0KHNlbGYuZmlsZSlcblxuICAgIGRlZiByZWFkX2ZpbGUoc2VsZik6XG4gICAgICAgIHdpdGggb3BlbihzZWxmLmZpbGUsICJyIikgYXMgZjpcbiAgICAgICAgICAgIHRleHQgPSBmLnJlYWRsaW5lcygpXG4gICAgICAgICAgICByZXR1cm4gdGV4dFxuXG4gICAgZGVmIHBhcnNlX2RvbWFpbihzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgZG9tYWluID0gW11cbiAgICAgICAgZm9yIGxpbmUgaW4gdGV4dDpcbiAgICAgICAgICAgIGlmIGxpbmUuc3RhcnRzd2l0aCgiZG9tYWluOiIpOlxuICAgICAgICAgICAgICAgIGRvbWFpbi5hcHBlbmQobGluZS5zcGxpdCgiOiIpWzFdLnN0cmlwKCkpXG4gICAgICAgIHJldHVybiBkb21haW5cblxuICAgIGRlZiBwYXJzZV9yb3V0ZShzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgcm91dGUgPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiB0ZXh0OlxuICAgICAgICAgICAgaWYgbGluZS5zdGFydHN3aXRoKCJyb3V0ZToiKTpcbiAgICAgICAgICAgICAgICByb3V0ZS5hcHBlbmQobGluZS5zcGxpdCgiOiIpWzFdLnN0cmlwKCkpXG4gICAgICAgIHJldHVybiByb3V0ZVxuXG4gICAgZGVmIHBhcnNlX25ldChzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgbmV0ID0gW11cbiAgICAgICAgZm9yIGxpbmUgaW4gdGV4dDpcbiAgICAgICAgICAgIGlmIGxpbmUuc3RhcnRzd2l0aCgibmV0OiIpOlxuICAgICAgICAgICAgICAgIG5ldC5hcHBlbmQobGluZS5zcGxpdCgiOiIpWzFdLnN0cmlwKCkpXG4gICAgICAgIHJldHVybiBuZXRcblxuICAgIGRlZiBwYXJzZV9hcyhzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgYXNfID0gW11cbiAgICAgICAgZm9yIGxpbmUgaW4gdGV4dDpcbiAgICAgICAgICAgIGlmIGxpbmUuc3RhcnRzd2l0aCgiYXM6Iik6XG4gICAgICAgICAgICAgICAgYXNfLmFwcGVuZChsaW5lLnNwbGl0KCI6IilbMV0uc3RyaXAoKSlcbiAgICAgICAgcmV0dXJuIGFzX1xuXG4gICAgZGVmIHBhcnNlX2Rlc2NyKHNlbGYpOlxuICAgICAgICB0ZXh0ID0gc2VsZi5yZWFkX2ZpbGUoKVxuICAgICAgICBkZXNjciA9IFtdXG4gICAgICAgIGZvciBsaW5lIGluIHRleHQ6XG4gICAgICAgICAgICBpZiBsaW5lLnN0YXJ0c3dpdGgoImRlc2NyOiIpOlxuICAgICAgICAgICAgICAgIGRlc2NyLmFwcGVuZChsaW5lLnNwbGl0KCI6IilbMV0uc3RyaXAoKSlcbiAgICAgICAgcmV0dXJuIGRlc2NyXG5cbiAgICBkZWYgcGFyc2Vfb3JpZ2luKHNlbGYpOlxuICAgICAgICB0ZXh0ID0gc2VsZi5yZWFkX2ZpbGUoKVxuICAgICAgICBvcmlnaW4gPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiB0ZXh0OlxuICAgICAgICAgICAgaWYgbGluZS5zdGFydHN3aXRoKCJvcmlnaW46Iik6XG4gICAgICAgICAgICAgICAgb3JpZ2luLmFwcGVuZChsaW5lLnNwbGl0KCI6IilbMV0uc3RyaXAoKSlcbiAgICAgICAgcmV0dXJuIG9yaWdpblxuXG4gICAgZGVmIHBhcnNlX21lbWJlcnMoc2VsZik6XG4gICAgICAgIHRleHQgPSBzZWxmLnJlYWRfZmlsZSgpXG4gICAgICAgIG1lbWJlcnMgPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiB0ZXh0OlxuICAgICAgICAgICAgaWYgbGluZS5zdGFydHN3aXRoKCJtZW1iZXJzOiIpOlxuICAgICAgICAgICAgICAgIG1lbWJlcnMuYXBwZW5kKGxpbmUuc3BsaXQoIjoiKVsxXS5zdHJpcCgpKVxuICAgICAgICByZXR1cm4gbWVtYmVyc1xuXG4gICAgZGVmIHBhcnNlX21udF9sb3dlcihzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgbW50X2xvd2VyID0gW11cbiAgICAgICAgZm9yIGxpbmUgaW4gdGV4dDpcbiAgICAgICAgICAgIGlmIGxpbmUuc3RhcnRzd2l0aCgibW50LWxvd2VyOiIpOlxuICAgICAgICAgICAgICAgIG1udF9sb3dlci5hcHBlbmQobGluZS5zcGxpdCgiOiIpWzFdLnN0cmlwKCkpXG4gICAgICAgIHJldHVybiBtbnRfbG93ZXJcblxuICAgIGRlZiBwYXJzZV9tbnRfcm91dGVzKHNlbGYpOlxuICAgICAgICB0ZXh0ID0gc2VsZi5yZWFkX2ZpbGUoKVxuICAgICAgICBtbnRfcm91dGVzID0gW11cbiAgICAgICAgZm9yIGxpbmUgaW4gdGV4dDpcbiAgICAgICAgICAgIGlmIGxpbmUuc3RhcnRzd2l0aCgibW50LXJvdXRlczoiKTpcbiAgICAgICAgICAgICAgICBtbnRfcm91dGVzLmFwcGVuZChsaW5lLnNwbGl0KCI6IilbMV0uc3RyaXAoKSlcbiAgICAgICAgcmV0dXJuIG1udF9yb3V0ZXNcblxuICAgIGRlZiBwYXJzZV9tbnRfYnkoc2VsZik6XG4gICAgICAgIHRleHQgPSBzZWxmLnJlYWRfZmlsZSgpXG4gICAgICAgIG1udF9ieSA9IFtdXG4gICAgICAgIGZvciBsaW5lIGluIHRleHQ6XG4gICAgICAgICAgICBpZiBsaW5lLnN0YXJ0c3dpdGgoIm1udC1ieToiKTpcbiAgICAgICAgICAgICAgICBtbnRfYnkuYXBwZW5kKGxpbmUuc3BsaXQoIjoiKVsxXS5zdHJpcCgpKVxuICAgICAgICByZXR1cm4gbW50X2J5XG5cbiAgICBkZWYgcGFyc2VfY2hhbmdlZChzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgY2hhbmdlZCA9IFtdXG4gICAgICAgIGZvciBsaW5lIGluIHRleHQ6XG4gICAgICAgICAgICBpZiBsaW5lLnN0YXJ0c3dpdGgoImNoYW5nZWQ6Iik6XG4gICAgICAgICAgICAgICAgY2hhbmdlZC5hcHBlbmQobGluZS5zcGxpdCgiOiIpWzFdLnN0cmlwKCkpXG4gICAgICAgIHJldHVybiBjaGFuZ2VkXG5cbiAgICBkZWYgcGFyc2Vfc291cmNlKHNlbGYpOlxuICAgICAgICB0ZXh0ID0gc2VsZi5yZWFkX2ZpbGUoKVxuICAgICAgICBzb3VyY2UgPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiB0ZXh0OlxuICAgICAgICAgICAgaWYgbGluZS5zdGFydHN3aXRoKCJzb3VyY2U6Iik6XG4gICAgICAgICAgICAgICAgc291cmNlLmFwcGVuZChsaW5lLnNwbGl0KCI6IilbMV0uc3RyaXAoKSlcbiAgICAgICAgcmV0dXJuIHNvdXJjZVxuXG4gICAgZGVmIHBhcnNlX3JlbWFya3Moc2VsZik6XG4gICAgICAgIHRleHQgPSBzZWxmLnJlYWRfZmlsZSgpXG4gICAgICAgIHJlbWFya3MgPSBbXVxuICAgICAgICBmb3IgbGluZSBpbiB0ZXh0OlxuICAgICAgICAgICAgaWYgbGluZS5zdGFydHN3aXRoKCJyZW1hcmtzOiIpOlxuICAgICAgICAgICAgICAgIHJlbWFya3MuYXBwZW5kKGxpbmUuc3BsaXQoIjoiKVsxXS5zdHJpcCgpKVxuICAgICAgICByZXR1cm4gcmVtYXJrc1xuXG4gICAgZGVmIHBhcnNlX25vdGlmeShzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgbm90aWZ5ID0gW11cbiAgICAgICAgZm9yIGxpbmUgaW4gdGV4dDpcbiAgICAgICAgICAgIGlmIGxpbmUuc3RhcnRzd2l0aCgibm90aWZ5OiIpOlxuICAgICAgICAgICAgICAgIG5vdGlmeS5hcHBlbmQobGluZS5zcGxpdCgiOiIpWzFdLnN0cmlwKCkpXG4gICAgICAgIHJldHVybiBub3RpZnlcblxuICAgIGRlZiBwYXJzZV9vcmcoc2VsZik6XG4gICAgICAgIHRleHQgPSBzZWxmLnJlYWRfZmlsZSgpXG4gICAgICAgIG9yZyA9IFtdXG4gICAgICAgIGZvciBsaW5lIGluIHRleHQ6XG4gICAgICAgICAgICBpZiBsaW5lLnN0YXJ0c3dpdGgoIm9yZzoiKTpcbiAgICAgICAgICAgICAgICBvcmcuYXBwZW5kKGxpbmUuc3BsaXQoIjoiKVsxXS5zdHJpcCgpKVxuICAgICAgICByZXR1cm4gb3JnXG5cbiAgICBkZWYgcGFyc2VfYWRtaW5fYyhzZWxmKTpcbiAgICAgICAgdGV4dCA9IHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgYWRtaW5fYyA9IFtdXG4gICAgICAgIGZvciBsaW5lIGluIHRleHQ6XG4gICAgICAgICAgICBpZiBsaW5lLnN0YXJ0c3dpdGgoImFkbWluLWM6Iik6XG4gICAgICAgICAgICAgICAgYWRtaW5fYy5hcHBlbmQobGluZS5zcGxpdCgiOiIpWzFdLnN0cmlwKCkpXG4gICAgICAgIHJldHVybiBhZG1pbl9jXG5cbicsICdcbmRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG5cbnNlbGYuZmlsZSA9IGZpbGVwYXRoXG5cblxuJywgJyAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBcbiAgICBkZWYgZ2V0X2ZpbGUoc2VsZik6XG4gICAgICAgIHJldHVybiBzZWxmLmZpbGVcbiAgICBcbiAgICBkZWYgZ2V0X2RvbWFpbihzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzFdXG4gICAgXG4gICAgZGVmIGdldF9wYXRoKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMF1cbiAgICBcbiAgICBkZWYgZ2V0X2ZpbGVuYW1lKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMl1cbiAgICBcbiAgICBkZWYgZ2V0X2Jhc2VuYW1lKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMl0uc3BsaXQoIi4iKVswXVxuICAgIFxuICAgIGRlZiBnZXRfZXh0ZW5zaW9uKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMl0uc3BsaXQoIi4iKVsxXVxuICAgIFxuICAgIGRlZiBnZXRfZmlsZV9wYXRoKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMF0gKyAiLyIgKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsxXVxuICAgIFxuICAgIGRlZiBnZXRfZmlsZW5hbWVfd2l0aF9wYXRoKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMF0gKyAiLyIgKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsxXSArICIvIiArIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzJdXG4gICAgXG4gICAgZGVmIGdldF9maWxlbmFtZV93aXRob3V0X2V4dGVuc2lvbihzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzJdLnNwbGl0KCIuIilbMF1cbiAgICBcbiAgICBkZWYgZ2V0X2ZpbGVuYW1lX3dpdGhvdXRfZXh0ZW5zaW9uX3dpdGhfcGF0aChzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzBdICsgIi8iICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMV0gKyAiLyIgKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsyXS5zcGxpdCgiLiIpWzBdXG4gICAgXG4gICAgZGVmIGdldF9maWxlbmFtZV93aXRob3V0X2V4dGVuc2lvbl93aXRoX2RvbWFpbihzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzFdICsgIi8iICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMl0uc3BsaXQoIi4iKVswXVxuICAgIFxuICAgIGRlZiBnZXRfZmlsZW5hbWVfd2l0aG91dF9leHRlbnNpb25fd2l0aF9wYXRoX2FuZF9kb21haW4oc2VsZik6XG4gICAgICAgIHJldHVybiBzZWxmLmZpbGUuc3BsaXQoIi8iKVswXSArICIvIiArIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzFdICsgIi8iICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMl0uc3BsaXQoIi4iKVswXVxuICAgIFxuICAgIGRlZiBnZXRfZmlsZW5hbWVfd2l0aG91dF9leHRlbnNpb25fd2l0aF9wYXRoX2FuZF9kb21haW5fY2FtZWxfY2FzZShzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzBdICsgIi8iICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMV0gKyAiLyIgKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsyXS5zcGxpdCgiLiIpWzBdLnRpdGxlKClcbiAgICBcbiAgICBkZWYgZ2V0X2ZpbGVuYW1lX3dpdGhfcGF0aF9hbmRfZG9tYWluX2NhbWVsX2Nhc2Uoc2VsZik6XG4gICAgICAgIHJldHVybiBzZWxmLmZpbGUuc3BsaXQoIi8iKVswXSArICIvIiArIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzFdICsgIi8iICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMl0udGl0bGUoKVxuICAgIFxuICAgIGRlZiBnZXRfZmlsZW5hbWVfY2FtZWxfY2FzZShzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzJdLnRpdGxlKClcbiAgICBcbiAgICBkZWYgZ2V0X2ZpbGVuYW1lX3dpdGhvdXRfZXh0ZW5zaW9uX2NhbWVsX2Nhc2Uoc2VsZik6XG4gICAgICAgIHJldHVybiBzZWxmLmZpbGUuc3BsaXQoIi8iKVsyXS5zcGxpdCgiLiIpWzBdLnRpdGxlKClcbiAgICBcbiAgICBkZWYgZ2V0X2ZpbGVfcGF0aF9jYW1lbF9jYXNlKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMF0gKyAiLyIgKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsxXS50aXRsZSgpXG4gICAgXG4gICAgZGVmIGdldF9maWxlX3BhdGhfY2FtZWxfY2FzZV93aXRob3V0X3NsYXNoKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMF0gKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsxXS50aXRsZSgpXG4gICAgXG4gICAgZGVmIGdldF9maWxlX3BhdGhfY2FtZWxfY2FzZV93aXRob3V0X3NsYXNoX3dpdGhvdXRfZG90KHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMF0gKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsxXS50aXRsZSgpICsgIi4iICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMl0uc3BsaXQoIi4iKVsxXVxuICAgIFxuICAgIGRlZiBnZXRfZmlsZV9wYXRoX2NhbWVsX2Nhc2Vfd2l0aG91dF9zbGFzaF93aXRob3V0X2RvdF93aXRob3V0X2V4dGVuc2lvbihzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzBdICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMV0udGl0bGUoKVxuICAgIFxuICAgIGRlZiBnZXRfZmlsZV9wYXRoX2NhbWVsX2Nhc2Vfd2l0aG91dF9zbGFzaF93aXRob3V0X2RvdF93aXRob3V0X2V4dGVuc2lvbl93aXRob3V0X2RvdChzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzBdICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMV0udGl0bGUoKSArICIuIlxuICAgIFxuICAgIGRlZiBnZXRfZmlsZV9wYXRoX3dpdGhvdXRfc2xhc2goc2VsZik6XG4gICAgICAgIHJldHVybiBzZWxmLmZpbGUuc3BsaXQoIi8iKVswXSArIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzFdXG4gICAgXG4gICAgZGVmIGdldF9maWxlX3BhdGhfd2l0aG91dF9zbGFzaF93aXRob3V0X2RvdChzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzBdICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMV0gKyAiLiIgKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsyXS5zcGxpdCgiLiIpWzFdXG4gICAgXG4gICAgZGVmIGdldF9maWxlX3BhdGhfd2l0aG91dF9zbGFzaF93aXRob3V0X2RvdF93aXRob3V0X2V4dGVuc2lvbihzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzBdICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMV1cbiAgICBcbiAgICBkZWYgZ2V0X2ZpbGVfcGF0aF93aXRob3V0X3NsYXNoX3dpdGhvdXRfZG90X3dpdGhvdXRfZXh0ZW5zaW9uX3dpdGhvdXRfZG90KHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMF0gKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsxXSArICIuIlxuICAgIFxuICAgIGRlZiBnZXRfZmlsZW5hbWVfd2l0aF9wYXRoX2FuZF9kb21haW5fY2FtZWxfY2FzZV93aXRob3V0X2RvdChzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzBdICsgIi8iICsgc2VsZi5maWxlLnNwbGl0KCIvIilbMV0gKyAiLyIgKyBzZWxmLmZpbGUuc3BsaXQoIi8iKVsyXS50aXRsZSgpICsgIi4iXG4gICAgXG4gICAgZGVmIGdldF9maWxlbmFtZV9jYW1lbF9jYXNlX3dpdGhvdXRfZG90KHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlLnNwbGl0KCIvIilbMl0udGl0bGUoKSArICIuIlxuICAgIFxuICAgIGRlZiBnZXRfZmlsZW5hbWVfd2l0aG91dF9leHRlbnNpb25fY2FtZWxfY2FzZV93aXRob3V0X2RvdChzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzJdLnNwbGl0KCIuIilbMF0udGl0bGUoKSArICIuIlxuICAgIFxuICAgIGRlZiBnZXRfZmlsZV9wYXRoX2NhbWVsX2Nhc2Vfd2l0aG91dF9kb3Qoc2VsZik6XG4gICAgICAgIHJldHVybiBzZWxmLmZpbGUuc3BsaXQoIi8iKVswXSArICIvIiArIHNlbGYuZmlsZS5zcGxpdCgiLyIpWzFdLnRpdGxlKCkgKyAiLiJcbiAgICBcbiAgICBkZWYgZ2V0X2ZpbGVfcGF0aF8="
Set-B64Value "G2" "W1RydWUsIFRydWUsIEZhbHNlLCBGYWxzZSwgVHJ1ZSwgVHJ1ZSwgRmFsc2UsIFRydWUsIFRydWUsIEZhbHNlXQ=="
$ws.Range("H2").Value = "KO"
$ws.Range("I2").Value = "OK"
$ws.Range("J2").Value = "OK"

# Row 3
Set-B64Value "B3" "RGVmaW5lIGEgY2xhc3MgbmFtZWQgUm91dGVEb21haW4gdGhhdCBoYXMgYW4gYXRyaWJ1dHRlIGNhbGxlZCBmaWxlLiBUaGUgYXR0cmlidXRlIGNhbGxlZCBmaWxlIGlzIGVxdWFsIHRvIGFuIGlucHV0IGdpdmVuIHRvIHRoZSBjbGFzcy4gVGhlIGlucHV0IGdpdmVuIHRvIHRoZSBjbGFzcyBpcyBjYWxsZWQgZmlsZXBhdGguIE9idGFpbiB0aGUgZXh0ZW5zaW9uIG9mIHRoZSBmaWxlLiBEZWZpbmUgYSBuZXcgYXR0cmlidXRlIGluIHRoZSBjbGFzcyBuYW1lZCBmaWxlX2VudGl0eSB0aGF0IGRlcGVuZGluZyBvZiB0aGUgZXh0ZW5zaW9uIG9mIHRoZSBmaWxlIGluc3RhbmNlcyBhIEZpbGVDU1YgY2xhc3Mgb3IgYSBGaWxlVFhUIGNsYXNzLiBEbyBub3QgdXNlIG9zIG9yIHN5cyBsaWJyYXJpZXM="
Set-B64Value "D3" "Y2xhc3MgUm91dGVEb21haW4oKQogICAgICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6CiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGgKICAgICAgICBleHQgPSBzZWxmLmZpbGUuc3BsaXQoJy4nKVstMV0KICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gZmlsZV9hY2NlcHRlZFtleHRdKHNlbGYuZmlsZSk="
$ws.Range("E3").Value = "test1"
Set-B64Value "F3" "WyIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IHNlbGYuZmlsZV9pbnN0YW5jZSgpXG4gICAgZGVmIGZpbGVfaW5zdGFuY2Uoc2VsZik6XG4gICAgICAgICMgZmlsZV9wYXRoID0gc2VsZi5maWxlXG4gICAgICAgIGV4dGVuc2lvbiA9IHNlbGYuZmlsZS5zcGxpdCgnLicpWy0xXVxuICAgICAgICBpZiBleHRlbnNpb24gPT0gJ2Nzdic6XG4gICAgICAgICAgICByZXR1cm4gRmlsZUNTVihzZWxmLmZpbGUpXG4gICAgICAgIGVsaWYgZXh0ZW5zaW9uID09ICd0eHQnOlxuICAgICAgICAgICAgcmV0dXJuIEZpbGVUWFQoc2VsZi5maWxlKVxuICAgICAgICBlbHNlOlxuICAgICAgICAgICAgcmFpc2UgVmFsdWVFcnJvcignRXh0ZW5zaW9uIG5vdCBhbGxvd2VkJylcblxuY2xhc3MgRmlsZUNTVigpOlxuICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZVxuICAgIGRlZiByZWFkKHNlbGYpOlxuICAgICAgICB3aXRoIG9wZW4oc2VsZi5maWxlLCAncicpIGFzIGY6XG4gICAgICAgICAgICByZXR1cm4gZi5yZWFkKClcblxuIiwgJyAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gTm9uZVxuICAgICAgICBzZWxmLmV4dGVuc2lvbiA9IHNlbGYuZmlsZS5zcGxpdChcJy5cJylbLTFdXG4gICAgICAgIGlmIHNlbGYuZXh0ZW5zaW9uID09IFwnY3N2XCc6XG4gICAgICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gRmlsZUNTVihzZWxmLmZpbGUpXG4gICAgICAgIGVsaWYgc2VsZi5leHRlbnNpb24gPT0gXCd0eHRcJzpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlVFhUKHNlbGYuZmlsZSlcbiAgICAgICAgZWxzZTpcbiAgICAgICAgICAgIHJhaXNlIFR5cGVFcnJvclxuXG5mcm9tIHJvdXRlX2RvbWFpbiBpbXBvcnQgUm91dGVEb21haW5cblxubG9nZ2luZy5iYXNpY0NvbmZpZyhsZXZlbD1sb2dnaW5nLklORk8pXG5cbmlmIF9fbmFtZV9fID09IFwnX19tYWluX19cJzpcbiAgICBwYXRoX2RhdGEgPSBcJ3JvdXRlLXZpZXdzLnJvdXRldmlld3Mub3JnX2JncGRhdGFfMjAxOS4wNi50eHRcJ1xuICAgIGRhdGEgPSBSb3V0ZURvbWFpbihwYXRoX2RhdGEpXG4gICAgbG9nZ2luZy5pbmZvKGRhdGEuZmlsZV9lbnRpdHkuYXNuX2FzX3ByZWZpeF9kaWN0KVxuMTU2MjE5MzYwMCwxNTYxOTMyMjYwLDAuMC4wLjAvMCwwLCIyODI4IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTI5NTYiLDAsMC4wLjAuMCwwLDEwMCwwLDAsMCwwLDAsMCwwXG4xNTYyMTkzNjAwLDE1NjE5MzIyNjAsMS4wLjAuMC8yNCw3NTQ1LCI3NTQ1IDEyOTkgMTInLCAiICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZSA9IGZpbGVwYXRoXG4gICAgICAgIHNlbGYuZmlsZV9leHRlbnNpb24gPSBzZWxmLmZpbGUuc3BsaXQoJy4nKVstMV1cbiAgICAgICAgaWYgc2VsZi5maWxlX2V4dGVuc2lvbiA9PSAnY3N2JzpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlQ1NWKHNlbGYuZmlsZSlcbiAgICAgICAgZWxpZiBzZWxmLmZpbGVfZXh0ZW5zaW9uID09ICd0eHQnOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVUWFQoc2VsZi5maWxlKVxuICAgICAgICBlbHNlOlxuICAgICAgICAgICAgcmFpc2UgRXhjZXB0aW9uKCdGaWxlIHR5cGUgbm90IHN1cHBvcnRlZCcpXG5cbmNsYXNzIEZpbGVDU1YoKTpcbiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGVfbmFtZSA9IHNlbGYuZmlsZS5zcGxpdCgnLycpWy0xXVxuICAgICAgICBcbiIsICIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IE5vbmVcbiAgICAgICAgc2VsZi5kZWZpbmVfZmlsZV9lbnRpdHkoKVxuICAgIFxuICAgIGRlZiBkZWZpbmVfZmlsZV9lbnRpdHkoc2VsZik6XG4gICAgICAgIGlmIHNlbGYuZmlsZS5lbmRzd2l0aCgnLmNzdicpOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVDU1Yoc2VsZi5maWxlKVxuICAgICAgICBlbGlmIHNlbGYuZmlsZS5lbmRzd2l0aCgnLnR4dCcpOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVUWFQoc2VsZi5maWxlKVxuICAgICAgICAgICAgXG4gICAgZGVmIGdldF9maWxlX2VudGl0eShzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZV9lbnRpdHlcbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUuY3N2JylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnR4dCcpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS54bHN4JylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlJylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnBkZicpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5kb2N4JylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnBuZycpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5qcGcnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUucmFyJylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLmd6aXAnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUuZ3onKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUuemlwJylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLmpzb24nKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUueG1sJylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnhscycpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS54bHNiJylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnhsc3gnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUueGx0JylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnhsdHgnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUub2RzJylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLmRiZicpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5kYicpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5zcWxpdGUnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUuc3FsaXRlMycpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5wcHR4JylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnBwdCcpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5wcHRtJylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnBwdHgnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUucHB0eCcpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5wcHR4JylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnBwdHgnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUucHB0eCcpXG5maWxlX2VudGl0eSA9IHJvdXRlX2RvbWFpbi5nZXRfZmlsZV9lbnRpdHkoKVxuZmlsZV9lbnRpdHkucmVhZF9maWxlKClcblxucm91dGVfZG9tYWluID0gUm91dGVEb21haW4oJy4uL2RhdGEvZmlsZS5wcHR4JylcbmZpbGVfZW50aXR5ID0gcm91dGVfZG9tYWluLmdldF9maWxlX2VudGl0eSgpXG5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKVxuXG5yb3V0ZV9kb21haW4gPSBSb3V0ZURvbWFpbignLi4vZGF0YS9maWxlLnBwdHgnKVxuZmlsZV9lbnRpdHkgPSByb3V0ZV9kb21haW4uZ2V0X2ZpbGVfZW50aXR5KClcbmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG5cbnJvdXRlX2RvbWFpbiA9IFJvdXRlRG9tYWluKCcuLi9kYXRhL2ZpbGUucHB0eCcpXG5maWxlX2VudGl0eSA9IiwgIiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmV4dGVuc2lvbiA9IHNlbGYuZmlsZS5zcGxpdCgnLicpWy0xXVxuICAgICAgICBpZiBzZWxmLmV4dGVuc2lvbiA9PSAnY3N2JzpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlQ1NWKHNlbGYuZmlsZSlcbiAgICAgICAgZWxzZTpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlVFhUKHNlbGYuZmlsZSlcblxuY2xhc3MgRmlsZUNTVigpOlxuICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZSA9IGZpbGVwYXRoXG5cbiIsICIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5leHRlbnNpb24gPSBzZWxmLmZpbGUuc3BsaXQoJy4nKVstMV1cbiAgICAgICAgaWYgc2VsZi5leHRlbnNpb24gPT0gJ2Nzdic6XG4gICAgICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gRmlsZUNTVihzZWxmLmZpbGUpXG4gICAgICAgIGVsc2U6XG4gICAgICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gRmlsZVRYVChzZWxmLmZpbGUpXG4gICAgICAgICAgICBcbiAgICBkZWYgcHJpbnRfZmlsZShzZWxmKTpcbiAgICAgICAgcHJpbnQoc2VsZi5maWxlX2VudGl0eS5yZWFkX2ZpbGUoKSlcblxuY2xhc3MgRmlsZShvYmplY3QpOlxuICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZXBhdGggPSBmaWxlcGF0aFxuICAgIFxuICAgIGRlZiByZWFkX2ZpbGUoc2VsZik6XG4gICAgICAgIHdpdGggb3BlbihzZWxmLmZpbGVwYXRoLCAncicpIGFzIGY6XG4gICAgICAgICAgICByZXR1cm4gZi5yZWFkKClcblxuIiwgIiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gTm9uZVxuICAgICAgICBzZWxmLmV4dGVuc2lvbiA9IHNlbGYuZmlsZS5zcGxpdCgnLicpWy0xXVxuXG4gICAgZGVmIF9fZ2V0YXR0cl9fKHNlbGYsIGF0dHIpOlxuICAgICAgICByZXR1cm4gZ2V0YXR0cihzZWxmLmZpbGVfZW50aXR5LCBhdHRyKVxuXG4gICAgZGVmIF9fcmVwcl9fKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlX2VudGl0eS5fX3JlcHJfXygpXG5cbiAgICBkZWYgX19zdHJfXyhzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZV9lbnRpdHkuX19zdHJfXygpXG5cbiAgICBkZWYgX19pdGVyX18oc2VsZik6XG4gICAgICAgIHJldHVybiBzZWxmLmZpbGVfZW50aXR5Ll9faXRlcl9fKClcblxuICAgIGRlZiBfX25leHRfXyhzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZV9lbnRpdHkuX19uZXh0X18oKVxuXG5jbGFzcyBGaWxlQ1NWKCk6XG4gICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5saW5lcyA9IE5vbmVcbiAgICAgICAgc2VsZi5yZWFkX2ZpbGUoKVxuXG4gICAgZGVmIHJlYWRfZmlsZShzZWxmKTpcbiAgICAgICAgd2l0aCBvcGVuKHNlbGYuZmlsZSwgJ3InKSBhcyBmaWxlOlxuICAgICAgICAgICAgc2VsZi5saW5lcyA9IGZpbGUucmVhZGxpbmVzKClcblxuICAgIGRlZiBfX3JlcHJfXyhzZWxmKTpcbiAgICAgICAgcmV0dXJuICdGaWxlQ1NWKHt9KScuZm9ybWF0KHNlbGYuZmlsZSlcblxuICAgIGRlZiBfX3N0cl9fKHNlbGYpOlxuICAgICAgICByZXR1cm4gJ0ZpbGVDU1Y6IHt9Jy5mb3JtYXQoc2VsZi5maWxlKVxuXG4gICAgZGVmIF9faXRlcl9fKHNlbGYpOlxuICAgICAgICByZXR1cm4gaXRlcihzZWxmLmxpbmVzKVxuXG4gICAgZGVmIF9fbmV4dF9fKHNlbGYpOlxuICAgICAgICByZXR1cm4gbmV4dChzZWxmLmxpbmVzKVxuXG4iLCAiICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZSA9IGZpbGVwYXRoXG4gICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBzZWxmLmNyZWF0ZV9lbnRpdHkoKVxuICAgIFxuICAgIGRlZiBjcmVhdGVfZW50aXR5KHNlbGYpOlxuICAgICAgICBpZiBzZWxmLmZpbGUuZW5kc3dpdGgoJy5jc3YnKTpcbiAgICAgICAgICAgIHJldHVybiBGaWxlQ1NWKHNlbGYuZmlsZSlcbiAgICAgICAgZWxpZiBzZWxmLmZpbGUuZW5kc3dpdGgoJy50eHQnKTpcbiAgICAgICAgICAgIHJldHVybiBGaWxlVFhUKHNlbGYuZmlsZSlcbiAgICAgICAgZWxzZTpcbiAgICAgICAgICAgIHJldHVybiBOb25lXG5cbmNsYXNzIEZpbGVDU1YoKTpcbiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmRhdGEgPSBzZWxmLnJlYWRfY3N2KClcbiAgICBcbiAgICBkZWYgcmVhZF9jc3Yoc2VsZik6XG4gICAgICAgIGRmID0gcGQucmVhZF9jc3Yoc2VsZi5maWxlKVxuICAgICAgICByZXR1cm4gZGZcblxuIiwgIiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gTm9uZVxuICAgICAgICBzZWxmLmV4dGVuc2lvbiA9IHNlbGYuZmlsZS5zcGxpdCgnLicpWy0xXVxuICAgICAgICBpZiBzZWxmLmV4dGVuc2lvbiA9PSAnY3N2JzpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlQ1NWKHNlbGYuZmlsZSlcbiAgICAgICAgZWxzZTpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlVFhUKHNlbGYuZmlsZSlcblxuY2xhc3MgRmlsZUNTVigpOlxuICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZVxuICAgICAgICBzZWxmLm51bWJlcl9vZl9saW5lcyA9IGxlbihsaXN0KG9wZW4oc2VsZi5maWxlKSkpXG5cbiIsICIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlcGF0aCA9IGZpbGVwYXRoXG4gICAgICAgIHNlbGYuZmlsZSA9IG9zLnBhdGguYmFzZW5hbWUoZmlsZXBhdGgpXG4gICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBzZWxmLmZpbGVfZW50aXR5KClcbiAgICBkZWYgZmlsZV9lbnRpdHkoc2VsZik6XG4gICAgICAgIGlmIHNlbGYuZmlsZS5zcGxpdCgnLicpWzFdID09ICdjc3YnOlxuICAgICAgICAgICAgZmlsZV9lbnRpdHkgPSBGaWxlQ1NWKHNlbGYuZmlsZXBhdGgpXG4gICAgICAgIGVsaWYgc2VsZi5maWxlLnNwbGl0KCcuJylbMV0gPT0gJ3R4dCc6XG4gICAgICAgICAgICBmaWxlX2VudGl0eSA9IEZpbGVUWFQoc2VsZi5maWxlcGF0aClcbiAgICAgICAgcmV0dXJuIGZpbGVfZW50aXR5XG4gICAgXG4gICAgXG5cblxuY2xhc3MgRmlsZUNTVigpOlxuICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZXBhdGggPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGUgPSBvcy5wYXRoLmJhc2VuYW1lKGZpbGVwYXRoKVxuICAgICAgICBzZWxmLmRmID0gc2VsZi5yZWFkX2NzdigpXG4gICAgZGVmIHJlYWRfY3N2KHNlbGYpOlxuICAgICAgICBkZiA9IHBkLnJlYWRfY3N2KHNlbGYuZmlsZXBhdGgpXG4gICAgICAgIHJldHVybiBkZi5oZWFkKClcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcbiAgICBcblxuIl0="
Set-B64Value "G3" "W1RydWUsIEZhbHNlLCBUcnVlLCBGYWxzZSwgVHJ1ZSwgVHJ1ZSwgRmFsc2UsIFRydWUsIFRydWUsIEZhbHNlXQ=="
$ws.Range("H3").Value = "KO"
$ws.Range("I3").Value = "OK"
$ws.Range("J3").Value = "OK"

# Row 4
Set-B64Value "B4" "RGVmaW5lIGEgY2xhc3MgbmFtZWQgUm91dGVEb21haW4gdGhhdCBoYXMgYW4gYXRyaWJ1dHRlIGNhbGxlZCBmaWxlLiBUaGUgYXR0cmlidXRlIGNhbGxlZCBmaWxlIGlzIGVxdWFsIHRvIGFuIGlucHV0IGdpdmVuIHRvIHRoZSBjbGFzcy4gVGhlIGlucHV0IGdpdmVuIHRvIHRoZSBjbGFzcyBpcyBjYWxsZWQgZmlsZXBhdGguIE9idGFpbiB0aGUgZXh0ZW5zaW9uIG9mIHRoZSBmaWxlLiBEZWZpbmUgYSBuZXcgYXR0cmlidXRlIGluIHRoZSBjbGFzcyBuYW1lZCBmaWxlX2VudGl0eSB0aGF0IGRlcGVuZGluZyBvZiB0aGUgZXh0ZW5zaW9uIG9mIHRoZSBmaWxlIGluc3RhbmNlcyBhIEZpbGVDU1YgY2xhc3Mgb3IgYSBGaWxlVFhUIGNsYXNzLiBVc2UgYSBmdW5jdGlvbiBjYWxsZWQgZ2V0X2xpbmVzIHRvIHJlYWQgZWFjaCBsaW5lIG9mIHRoZSBmaWxlX2VudGl0eSBhbmQgc3RvcmUgaXQgaW4gYSB2YXJpYWJsZSBjYWxsZWQgZGF0YS4gIERlZmluZSBhIG5ldyBhdHRyaWJ1dGUgaW4gdGhlIGNsYXNzIG5hbWVkIGdyYXBoIHRoYXQgdXNpbmcgdGhlIHZhcmlhYmxlIGRhdGEgaW5zdGFuY2VzIGEgR3JhcGggY2xhc3MuIERvIG5vdCB1c2Ugb3Mgb3Igc3lzIGxpYnJhcmllcw=="
Set-B64Value "D4" "Y2xhc3MgUm91dGVEb21haW4oKQogICAgICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6CiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGgKICAgICAgICBleHQgPSBzZWxmLmZpbGUuc3BsaXQoJy4nKVstMV0KICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gZmlsZV9hY2NlcHRlZFtleHRdKHNlbGYuZmlsZSkKICAgICAgICBkYXRhID0gc2VsZi5maWxlX2VudGl0eS5nZXRfbGluZXMoKQogICAgICAgIHNlbGYuZ3JhcGggPSBHcmFwaChkYXRhKQ=="
$ws.Range("E4").Value = "test2"
Set-B64Value "F4" "WyIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IE5vbmVcbiAgICAgICAgc2VsZi5ncmFwaCA9IE5vbmVcbiAgICAgICAgc2VsZi5kYXRhID0gTm9uZVxuICAgICAgICBzZWxmLmdldF9saW5lcygpXG4gICAgICAgIHNlbGYuZ2V0X2ZpbGVfZW50aXR5KClcbiAgICAgICAgc2VsZi5nZXRfZ3JhcGgoKVxuICAgICAgICBcbiAgICBkZWYgZ2V0X2xpbmVzKHNlbGYpOlxuICAgICAgICB3aXRoIG9wZW4oc2VsZi5maWxlKSBhcyBmOlxuICAgICAgICAgICAgc2VsZi5kYXRhID0gZi5yZWFkbGluZXMoKVxuICAgICAgICAgICAgXG4gICAgZGVmIGdldF9maWxlX2VudGl0eShzZWxmKTpcbiAgICAgICAgaWYgc2VsZi5maWxlLmVuZHN3aXRoKCcuY3N2Jyk6XG4gICAgICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gRmlsZUNTVihzZWxmLmZpbGUpXG4gICAgICAgIGVsaWYgc2VsZi5maWxlLmVuZHN3aXRoKCcudHh0Jyk6XG4gICAgICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gRmlsZVRYVChzZWxmLmZpbGUpXG4gICAgICAgIGVsc2U6XG4gICAgICAgICAgICBwcmludCgnRmlsZSBub3Qgc3VwcG9ydGVkJylcbiAgICAgICAgICAgIFxuICAgIGRlZiBnZXRfZ3JhcGgoc2VsZik6XG4gICAgICAgIHNlbGYuZ3JhcGggPSBHcmFwaChzZWxmLmRhdGEpXG5cblxuY2xhc3MgRmlsZUNTVigpOlxuICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZSA9IGZpbGVwYXRoXG4gICAgICAgIHNlbGYuZGF0YSA9IFtdXG4gICAgICAgIHNlbGYuZ2V0X2RhdGEoKVxuICAgICAgICBcbiAgICBkZWYgZ2V0X2RhdGEoc2VsZik6XG4gICAgICAgIHdpdGggb3BlbihzZWxmLmZpbGUpIGFzIGY6XG4gICAgICAgICAgICBzZWxmLmRhdGEgPSBmLnJlYWRsaW5lcygpXG5cblxuIiwgIlx0ZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcblx0XHRzZWxmLmZpbGUgPSBmaWxlcGF0aFxuXHRcdHRyeTpcblx0XHRcdHNlbGYuZmlsZV9lbnRpdHkgPSAoRmlsZUNTVihzZWxmLmZpbGUpIGlmIHNlbGYuZmlsZS5zcGxpdCgnLicpWzFdID09ICdjc3YnIGVsc2UgRmlsZVRYVChzZWxmLmZpbGUpKVxuXHRcdFx0c2VsZi5kYXRhID0gc2VsZi5maWxlX2VudGl0eS5nZXRfbGluZXMoKVxuXHRcdFx0c2VsZi5ncmFwaCA9IEdyYXBoKHNlbGYuZGF0YSlcblx0XHRleGNlcHQgRXhjZXB0aW9uIGFzIGU6XG5cdFx0XHRwcmludChlKVxuXG5cbmNsYXNzIEdyYXBoKCk6XG5cdGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlX2RhdGEpOlxuXHRcdHNlbGYuZGF0YSA9IGZpbGVfZGF0YVxuXHRcdHNlbGYuZyA9IHt9XG5cdFx0Zm9yIGxpbmUgaW4gc2VsZi5kYXRhOlxuXHRcdFx0c2VsZi5nW2xpbmVbMF1dID0gW2xpbmVbMTpdXVxuXHRcdHNlbGYudmVydGljZXMgPSBzZWxmLmdldF92ZXJ0aWNlcygpXG5cdFx0c2VsZi5lZGdlcyA9IHNlbGYuZ2V0X2VkZ2VzKClcblxuXHRkZWYgZ2V0X3ZlcnRpY2VzKHNlbGYpOlxuXHRcdHJldHVybiBsaXN0KHNlbGYuZy5rZXlzKCkpXG5cblx0ZGVmIGdldF9lZGdlcyhzZWxmKTpcblx0XHRyZXR1cm4gbGlzdChzZWxmLmcudmFsdWVzKCkpXG5cblx0ZGVmIGdldF9wYXRoKHNlbGYsIG9yaWdpbik6XG5cdFx0c2VsZi5vcmlnaW4gPSBvcmlnaW5cblx0XHRyZXR1cm4gc2VsZi5nZXRfcGF0aF9hdXgoc2VsZi5vcmlnaW4sIFtdKVxuXG5cdGRlZiBnZXRfY29zdChzZWxmLCBvcmlnaW4pOlxuXHRcdHNlbGYub3JpZ2luID0gb3JpZ2luXG5cdFx0cmV0dXJuIHNlbGYuZ2V0X2Nvc3RfYXV4KHNlbGYub3JpZ2luLCAwKVxuXG5cdGRlZiBnZXRfdGltZShzZWxmLCBvcmlnaW4pOlxuXHRcdHNlbGYub3JpZ2luID0gb3JpZ2luXG5cdFx0cmV0dXJuIHNlbGYuZ2V0X3RpbWVfYXV4KHNlbGYub3JpZ2luLCAwKVxuXG5cdGRlZiBnZXRfcGF0aF9hdXgoc2VsZiwgb3JpZ2luLCB2aXNpdGVkKTpcblx0XHRzZWxmLm9yaWdpbiA9IG9yaWdpblxuXHRcdHNlbGYudmlzaXRlZCA9IHZpc2l0ZWRcblx0XHRzZWxmLnZpc2l0ZWQuYXBwZW5kKHNlbGYub3JpZ2luKVxuXHRcdGZvciBlZGdlIGluIHNlbGYuZWRnZXM6XG5cdFx0XHRpZiBzZWxmLm9yaWdpbiBpbiBlZGdlWzBdOlxuXHRcdFx0XHRmb3IgdmVydGV4IGluIGVkZ2VbMF06XG5cdFx0XHRcdFx0aWYgdmVydGV4IG5vdCBpbiBzZWxmLnZpc2l0ZWQ6XG5cdFx0XHRcdFx0XHRzZWxmLmdldF9wYXRoX2F1eCh2ZXJ0ZXgsIHNlbGYudmlzaXRlZClcblx0XHRyZXR1cm4gc2VsZi52aXNpdGVkXG5cblx0ZGVmIGdldF9jb3N0X2F1eChzZWxmLCBvcmlnaW4sIGNvc3QpOlxuXHRcdHNlbGYub3JpZ2luID0gb3JpZ2luXG5cdFx0c2VsZi5jb3N0ID0gY29zdFxuXHRcdGZvciBlZGdlIGluIHNlbGYuZWRnZXM6XG5cdFx0XHRpZiBzZWxmLm9yaWdpbiBpbiBlZGdlWzBdOlxuXHRcdFx0XHRmb3IgdmVydGV4IGluIGVkZ2VbMF06XG5cdFx0XHRcdFx0aWYgdmVydGV4IGluIHNlbGYudmlzaXRlZDpcblx0XHRcdFx0XHRcdHNlbGYuY29zdCArPSBpbnQoZWRnZVswXVstMV0pXG5cdFx0XHRcdFx0XHRzZWxmLmdldF9jb3N0X2F1eCh2ZXJ0ZXgsIHNlbGYuY29zdClcblx0XHRyZXR1cm4gc2VsZi5jb3N0XG5cblx0ZGVmIGdldF90aW1lX2F1eChzZWxmLCBvcmlnaW4sIHRpbWUpOlxuXHRcdHNlbGYub3JpZ2luID0gb3JpZ2luXG5cdFx0c2VsZi50aW1lID0gdGltZVxuXHRcdGZvciBlZGdlIGluIHNlbGYuZWRnZXM6XG5cdFx0XHRpZiBzZWxmLm9yaWdpbiBpbiBlZGdlWzBdOlxuXHRcdFx0XHRmb3IgdmVydGV4IGluIGVkZ2VbMF06XG5cdFx0XHRcdFx0aWYgdmVydGV4IGluIHNlbGYudmlzaXRlZDpcblx0XHRcdFx0XHRcdHNlbGYudGltZSArPSBpbnQoZWRnZVswXVstMl0pXG5cdFx0XHRcdFx0XHRzZWxmLmdldF90aW1lX2F1eCh2ZXJ0ZXgsIHNlbGYudGltZSlcblx0XHRyZXR1cm4gc2VsZi50aW1lXG5cblxuIiwgIiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gTm9uZVxuICAgICAgICBzZWxmLmdyYXBoID0gTm9uZVxuICAgICAgICBzZWxmLmdldF9saW5lcygpXG4gICAgICAgIHNlbGYuY3JlYXRlX2dyYXBoKClcbiAgICBcbiAgICBkZWYgZ2V0X2V4dGVuc2lvbihzZWxmKTpcbiAgICAgICAgZXh0ZW5zaW9uID0gc2VsZi5maWxlLnNwbGl0KCcuJylbLTFdXG4gICAgICAgIHJldHVybiBleHRlbnNpb25cbiAgICBcbiAgICBkZWYgY3JlYXRlX2ZpbGVfZW50aXR5KHNlbGYpOlxuICAgICAgICBpZiBzZWxmLmdldF9leHRlbnNpb24oKSA9PSAnY3N2JzpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlQ1NWKHNlbGYuZmlsZSlcbiAgICAgICAgZWxpZiBzZWxmLmdldF9leHRlbnNpb24oKSA9PSAndHh0JzpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlVFhUKHNlbGYuZmlsZSlcbiAgICBcbiAgICBkZWYgZ2V0X2xpbmVzKHNlbGYpOlxuICAgICAgICBzZWxmLmNyZWF0ZV9maWxlX2VudGl0eSgpXG4gICAgICAgIGRhdGEgPSBzZWxmLmZpbGVfZW50aXR5LnJlYWRfZmlsZSgpXG4gICAgICAgIHJldHVybiBkYXRhXG4gICAgXG4gICAgZGVmIGNyZWF0ZV9ncmFwaChzZWxmKTpcbiAgICAgICAgZGF0YSA9IHNlbGYuZ2V0X2xpbmVzKClcbiAgICAgICAgc2VsZi5ncmFwaCA9IEdyYXBoKGRhdGEpXG5cblxuY2xhc3MgRmlsZSgpOlxuICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZSA9IGZpbGVwYXRoXG4gICAgICAgIHNlbGYucmVhZF9maWxlKClcbiAgICAgICAgXG4gICAgZGVmIHJlYWRfZmlsZShzZWxmKTpcbiAgICAgICAgd2l0aCBvcGVuKHNlbGYuZmlsZSwgJ3InKSBhcyBmaWxlOlxuICAgICAgICAgICAgZGF0YSA9IGZpbGUucmVhZGxpbmVzKClcbiAgICAgICAgZGF0YSA9IFtsaW5lLnN0cmlwKCdcbicpIGZvciBsaW5lIGluIGRhdGFdXG4gICAgICAgIHJldHVybiBkYXRhXG4gICAgXG4iLCAiICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZSA9IGZpbGVwYXRoXG4gICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBOb25lXG4gICAgICAgIHNlbGYuZ3JhcGggPSBOb25lXG4gICAgICAgIHNlbGYuZ2V0X2xpbmVzKClcbiAgICAgICAgc2VsZi5jcmVhdGVfZ3JhcGgoKVxuICAgIFxuICAgIGRlZiBnZXRfbGluZXMoc2VsZik6XG4gICAgICAgIGZpbGVfZXh0ZW5zaW9uID0gc2VsZi5maWxlLnNwbGl0KCcuJylbLTFdXG4gICAgICAgIGlmIGZpbGVfZXh0ZW5zaW9uID09ICdjc3YnOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVDU1Yoc2VsZi5maWxlKVxuICAgICAgICBlbGlmIGZpbGVfZXh0ZW5zaW9uID09ICd0eHQnOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVUWFQoc2VsZi5maWxlKVxuICAgICAgICBlbHNlOlxuICAgICAgICAgICAgcHJpbnQoJ0Vycm9yJylcbiAgICBcbiAgICBkZWYgY3JlYXRlX2dyYXBoKHNlbGYpOlxuICAgICAgICBzZWxmLmdyYXBoID0gR3JhcGgoc2VsZi5maWxlX2VudGl0eS5kYXRhKVxuXG5cbmNsYXNzIEZpbGVDU1YoRmlsZSk6XG4gICAgZGVmIHJlYWQoc2VsZik6XG4gICAgICAgIGRhdGEgPSBbXVxuICAgICAgICB3aXRoIG9wZW4oc2VsZi5maWxlLCAncicpIGFzIGY6XG4gICAgICAgICAgICByZWFkZXIgPSBjc3YucmVhZGVyKGYpXG4gICAgICAgICAgICBmb3Igcm93IGluIHJlYWRlcjpcbiAgICAgICAgICAgICAgICBkYXRhLmFwcGVuZChyb3cpXG4gICAgICAgIHNlbGYuZGF0YSA9IGRhdGFcblxuXG4iLCAiICAgIGRlZiBfX2luaXRfXyhzZWxmLCBmaWxlcGF0aCk6XG4gICAgICAgIHNlbGYuZmlsZSA9IGZpbGVwYXRoXG4gICAgICAgIHNlbGYuZXh0ZW5zaW9uID0gc2VsZi5maWxlLnNwbGl0KCcuJylbLTFdXG4gICAgICAgIGlmIHNlbGYuZXh0ZW5zaW9uID09ICdjc3YnOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVDU1Yoc2VsZi5maWxlKVxuICAgICAgICBlbGlmIHNlbGYuZXh0ZW5zaW9uID09ICd0eHQnOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVUWFQoc2VsZi5maWxlKVxuICAgICAgICBzZWxmLmRhdGEgPSBzZWxmLmZpbGVfZW50aXR5LmdldF9saW5lcygpXG4gICAgICAgIHNlbGYuZ3JhcGggPSBHcmFwaChzZWxmLmRhdGEpXG5cblxuY2xhc3MgRmlsZUNTVihGaWxlKTpcbiAgICBkZWYgZ2V0X2xpbmVzKHNlbGYpOlxuICAgICAgICBkYXRhID0gW11cbiAgICAgICAgd2l0aCBvcGVuKHNlbGYuZmlsZSwgJ3InKSBhcyBmOlxuICAgICAgICAgICAgcmVhZGVyID0gY3N2LnJlYWRlcihmKVxuICAgICAgICAgICAgZm9yIHJvdyBpbiByZWFkZXI6XG4gICAgICAgICAgICAgICAgZGF0YS5hcHBlbmQocm93KVxuICAgICAgICByZXR1cm4gZGF0YVxuXG5cbiIsICIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IHNlbGYuZ2V0X2ZpbGVfZW50aXR5KClcbiAgICAgICAgc2VsZi5ncmFwaCA9IHNlbGYuZ2V0X2dyYXBoKClcblxuICAgIGRlZiBnZXRfZmlsZV9lbnRpdHkoc2VsZik6XG4gICAgICAgIGV4dGVuc2lvbiA9IHNlbGYuZmlsZS5zcGxpdCgnLicpWy0xXVxuICAgICAgICBpZiBleHRlbnNpb24gPT0gJ2Nzdic6XG4gICAgICAgICAgICBmaWxlX2VudGl0eSA9IEZpbGVDU1Yoc2VsZi5maWxlKVxuICAgICAgICBlbGlmIGV4dGVuc2lvbiA9PSAndHh0JzpcbiAgICAgICAgICAgIGZpbGVfZW50aXR5ID0gRmlsZVRYVChzZWxmLmZpbGUpXG4gICAgICAgIHJldHVybiBmaWxlX2VudGl0eVxuXG4gICAgZGVmIGdldF9ncmFwaChzZWxmKTpcbiAgICAgICAgZGF0YSA9IHNlbGYuZmlsZV9lbnRpdHkuZ2V0X2xpbmVzKClcbiAgICAgICAgZ3JhcGggPSBHcmFwaChkYXRhKVxuICAgICAgICByZXR1cm4gZ3JhcGhcblxuXG5jbGFzcyBGaWxlQ1NWKCk6XG4gICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcblxuICAgIGRlZiBnZXRfbGluZXMoc2VsZik6XG4gICAgICAgIHdpdGggb3BlbihzZWxmLmZpbGUsICdyJykgYXMgZjpcbiAgICAgICAgICAgIGRhdGEgPSBmLnJlYWRsaW5lcygpXG4gICAgICAgIHJldHVybiBkYXRhXG5cblxuIiwgIiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGVfZW50aXR5ID0gc2VsZi5nZXRfZmlsZV9lbnRpdHkoKVxuICAgICAgICBzZWxmLmRhdGEgPSBzZWxmLmdldF9saW5lcygpXG4gICAgICAgIHNlbGYuZ3JhcGggPSBzZWxmLmNyZWF0ZV9ncmFwaCgpXG5cbiAgICBkZWYgZ2V0X2V4dGVuc2lvbihzZWxmKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZmlsZS5zcGxpdCgnLicpWy0xXVxuXG4gICAgZGVmIGdldF9maWxlX2VudGl0eShzZWxmKTpcbiAgICAgICAgaWYgc2VsZi5nZXRfZXh0ZW5zaW9uKCkgPT0gJ2Nzdic6XG4gICAgICAgICAgICByZXR1cm4gRmlsZUNTVihzZWxmLmZpbGUpXG4gICAgICAgIGVsc2U6XG4gICAgICAgICAgICByZXR1cm4gRmlsZVRYVChzZWxmLmZpbGUpXG5cbiAgICBkZWYgZ2V0X2xpbmVzKHNlbGYpOlxuICAgICAgICByZXR1cm4gc2VsZi5maWxlX2VudGl0eS5nZXRfbGluZXMoKVxuXG4gICAgZGVmIGNyZWF0ZV9ncmFwaChzZWxmKTpcbiAgICAgICAgcmV0dXJuIEdyYXBoKHNlbGYuZGF0YSlcblxuXG5jbGFzcyBGaWxlQ1NWKCk6XG4gICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5saW5lcyA9IHNlbGYuZ2V0X2xpbmVzKClcblxuICAgIGRlZiBnZXRfbGluZXMoc2VsZik6XG4gICAgICAgIHdpdGggb3BlbihzZWxmLmZpbGUsICdyJykgYXMgZmlsZTpcbiAgICAgICAgICAgIHJldHVybiBbbGluZS5zdHJpcCgpIGZvciBsaW5lIGluIGZpbGUucmVhZGxpbmVzKCldXG5cblxuIiwgIlx0ZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcblx0XHRzZWxmLmZpbGUgPSBmaWxlcGF0aFxuXHRcdHNlbGYuZmlsZV9lbnRpdHkgPSBzZWxmLmZpbGVfZW50aXR5KClcblx0XHRzZWxmLmRhdGEgPSBzZWxmLmdldF9saW5lcygpXG5cdFx0c2VsZi5ncmFwaCA9IEdyYXBoKHNlbGYuZGF0YSlcblxuXHRkZWYgZmlsZV9lbnRpdHkoc2VsZik6XG5cdFx0aWYgc2VsZi5maWxlLmVuZHN3aXRoKCcuY3N2Jyk6XG5cdFx0XHRmaWxlX2VudGl0eSA9IEZpbGVDU1Yoc2VsZi5maWxlKVxuXHRcdGVsaWYgc2VsZi5maWxlLmVuZHN3aXRoKCcudHh0Jyk6XG5cdFx0XHRmaWxlX2VudGl0eSA9IEZpbGVUWFQoc2VsZi5maWxlKVxuXHRcdHJldHVybiBmaWxlX2VudGl0eVxuXG5cdGRlZiBnZXRfbGluZXMoc2VsZik6XG5cdFx0cmV0dXJuIHNlbGYuZmlsZV9lbnRpdHkucmVhZF9saW5lcygpXG5cblx0ZGVmIGdldF9ncmFwaChzZWxmKTpcblx0XHRyZXR1cm4gc2VsZi5ncmFwaFxuXG5cbmNsYXNzIEZpbGVDU1YoKTpcblx0ZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcblx0XHRzZWxmLmZpbGUgPSBmaWxlcGF0aFxuXG5cdGRlZiByZWFkX2xpbmVzKHNlbGYpOlxuXHRcdGRhdGEgPSBbXVxuXHRcdHdpdGggb3BlbihzZWxmLmZpbGUsICdyJykgYXMgZjpcblx0XHRcdHJlYWRlciA9IGNzdi5yZWFkZXIoZilcblx0XHRcdGZvciBsaW5lIGluIHJlYWRlcjpcblx0XHRcdFx0ZGF0YS5hcHBlbmQobGluZSlcblx0XHRyZXR1cm4gZGF0YVxuXG5cbiIsICIgICAgZGVmIF9faW5pdF9fKHNlbGYsIGZpbGVwYXRoKTpcbiAgICAgICAgc2VsZi5maWxlID0gZmlsZXBhdGhcbiAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IHNlbGYuZmlsZV9tYW5hZ2VyKClcbiAgICAgICAgc2VsZi5kYXRhID0gc2VsZi5nZXRfbGluZXMoKVxuICAgICAgICBzZWxmLmdyYXBoID0gR3JhcGgoc2VsZi5kYXRhKVxuXG4gICAgZGVmIGZpbGVfbWFuYWdlcihzZWxmKTpcbiAgICAgICAgZXh0ZW5zaW9uID0gc2VsZi5maWxlLnNwbGl0KCcuJylbLTFdXG4gICAgICAgIGlmIGV4dGVuc2lvbiA9PSAnY3N2JzpcbiAgICAgICAgICAgIHJldHVybiBGaWxlQ1NWKHNlbGYuZmlsZSlcbiAgICAgICAgZWxpZiBleHRlbnNpb24gPT0gJ3R4dCc6XG4gICAgICAgICAgICByZXR1cm4gRmlsZVRYVChzZWxmLmZpbGUpXG5cbiAgICBkZWYgZ2V0X2xpbmVzKHNlbGYpOlxuICAgICAgICBkYXRhID0gc2VsZi5maWxlX2VudGl0eS5yZWFkKClcbiAgICAgICAgcmV0dXJuIGRhdGFcblxuXG5cblxuY2xhc3MgR3JhcGgoKTpcbiAgICBkZWYgX19pbml0X18oc2VsZiwgZGF0YSk6XG4gICAgICAgIHNlbGYuZWRnZXMgPSBzZWxmLmdyYXBoX2VkZ2VzKGRhdGEpXG5cbiAgICBkZWYgZ3JhcGhfZWRnZXMoc2VsZiwgZGF0YSk6XG4gICAgICAgIGVkZ2VzID0gW11cbiAgICAgICAgZm9yIGl0ZW0gaW4gZGF0YTpcbiAgICAgICAgICAgIGVkZ2VzLmFwcGVuZCgoaXRlbVswXSwgaXRlbVsxXSwgaXRlbVsyXSkpXG4gICAgICAgIHJldHVybiBlZGdlc1xuXG5cblxuIiwgIiAgICBkZWYgX19pbml0X18oc2VsZiwgZmlsZXBhdGgpOlxuICAgICAgICBzZWxmLmZpbGUgPSBmaWxlcGF0aFxuICAgICAgICBzZWxmLmZpbGVfZXh0ZW5zaW9uID0gZmlsZXBhdGguc3BsaXQoJy4nKVstMV1cbiAgICAgICAgaWYgc2VsZi5maWxlX2V4dGVuc2lvbiA9PSAnY3N2JzpcbiAgICAgICAgICAgIHNlbGYuZmlsZV9lbnRpdHkgPSBGaWxlQ1NWKGZpbGVwYXRoKVxuICAgICAgICBlbHNlOlxuICAgICAgICAgICAgc2VsZi5maWxlX2VudGl0eSA9IEZpbGVUWFQoZmlsZXBhdGgpXG4gICAgICAgIHNlbGYuZGF0YSA9IHNlbGYuZmlsZV9lbnRpdHkuZ2V0X2xpbmVzKClcbiAgICAgICAgc2VsZi5ncmFwaCA9IEdyYXBoKHNlbGYuZGF0YSlcblxuICAgIGRlZiBmaW5kX3Nob3J0ZXN0X3BhdGhfYmV0d2Vlbl90d29fdmVydGljZXMoc2VsZiwgb3JpZ2luLCBkZXN0aW55KTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZ3JhcGguZmluZF9zaG9ydGVzdF9wYXRoX2JldHdlZW5fdHdvX3ZlcnRpY2VzKG9yaWdpbiwgZGVzdGlueSlcblxuXG5cblxuXG4gICAgZGVmIGZpbmRfZGlzdGFuY2VfYmV0d2Vlbl90d29fdmVydGljZXMoc2VsZiwgb3JpZ2luLCBkZXN0aW55KTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZ3JhcGguZmluZF9kaXN0YW5jZV9iZXR3ZWVuX3R3b192ZXJ0aWNlcyhvcmlnaW4sIGRlc3RpbnkpXG5cbiAgICBkZWYgZmluZF9udW1iZXJfb2ZfdHJpcHNfYmV0d2Vlbl90d29fdmVydGljZXMoc2VsZiwgb3JpZ2luLCBkZXN0aW55LCBtYXhpbXVtX3N0b3BzKTpcbiAgICAgICAgcmV0dXJuIHNlbGYuZ3JhcGguZmluZF9udW1iZXJfb2ZfdHJpcHNfYmV0d2Vlbl90d29fdmVydGljZXMob3JpZ2luLCBkZXN0aW55LCBtYXhpbXVtX3N0b3BzKVxuXG4gICAgZGVmIGZpbmRfbnVtYmVyX29mX3RyaXBzX2JldHdlZW5fdHdvX3ZlcnRpY2VzX3dpdGhfZXhhY3Rfc3RvcHMoc2VsZiwgb3JpZ2luLCBkZXN0aW55LCBleGFjdF9zdG9wcyk6XG4gICAgICAgIHJldHVybiBzZWxmLmdyYXBoLmZpbmRfbnVtYmVyX29mX3RyaXBzX2JldHdlZW5fdHdvX3ZlcnRpY2VzX3dpdGhfZXhhY3Rfc3RvcHMob3JpZ2luLCBkZXN0aW55LCBleGFjdF9zdG9wcylcblxuICAgIGRlZiBmaW5kX2xlbmd0aF9vZl9zaG9ydGVzdF9yb3V0ZV9iZXR3ZWVuX3R3b192ZXJ0aWNlcyhzZWxmLCBvcmlnaW4sIGRlc3RpbnksIG1heGltdW1fZGlzdGFuY2UpOlxuICAgICAgICByZXR1cm4gc2VsZi5ncmFwaC5maW5kX2xlbmd0aF9vZl9zaG9ydGVzdF9yb3V0ZV9iZXR3ZWVuX3R3b192ZXJ0aWNlcyhvcmlnaW4sIGRlc3RpbnksIG1heGltdW1fZGlzdGFuY2UpIl0="
Set-B64Value "G4" "W0ZhbHNlLCBUcnVlLCBGYWxzZSwgRmFsc2UsIFRydWUsIFRydWUsIFRydWUsIEZhbHNlLCBGYWxzZSwgVHJ1ZV0="
$ws.Range("H4").Value = "KO"
$ws.Range("I4").Value = "OK"
$ws.Range("J4").Value = "OK"
